# "completed creation of 1 min sheets with formatting"
# Refresh the CLOSE column (column D, rows 22-30) on Sheet1: the minute bar
# closed out, so every row's CLOSE value shifts down into the next row, and
# the newly-opened row 22 picks up the LTP value that was already showing
# in E22.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")
$ws.Activate()

$ws.Range("D30").Value = 8678.9
$ws.Range("D29").Value = 120.8
$ws.Range("D28").Value = 257.35000000000002
$ws.Range("D27").Value = 653.5
$ws.Range("D26").Value = 947.15
$ws.Range("D25").Value = 672.5
$ws.Range("D24").Value = 581.29999999999995
$ws.Range("D23").Value = 2315.1
$ws.Range("D22").Value = 245.7

# Recalculate so the refreshed figures are fully committed.
$excel.CalculateFullRebuild()

# Leave the selection where the author finished working.
$ws.Range("I10").Select()
